$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos:" block (rows 25-27) lists three course-requirement lines in
# both column B and column C (identical text, mirrored). The edit re-orders
# them so the "LOM3246 ... (Indicação de Conjunto)" line moves from the first
# slot to the last slot, while the other two lines shift up:
#
#   before                                                  after
#   row25: LOM3246 (Indicação de Conjunto)        ->  row25: LOB1021 (Requisito)
#   row26: LOB1021 (Requisito)                    ->  row26: LOM3016 (Requisito)
#   row27: LOM3016 (Requisito)                    ->  row27: LOM3246 (Indicação de Conjunto)

$linhaLOB1021 = "LOB1021 -  Física IV  (Requisito)`n"
$linhaLOM3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$linhaLOM3246 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"

$ws.Range("B25:C25").Value = $linhaLOB1021
$ws.Range("B26:C26").Value = $linhaLOM3016
$ws.Range("B27:C27").Value = $linhaLOM3246
